# Update column F ("dSF") values in several rows to match repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    13 = -2
    15 = 7
    24 = 3
    27 = -1
    45 = 1
    47 = 2
    52 = -5
    53 = -4
    57 = 7
    61 = -4
    64 = -3
    66 = 2
    74 = -4
    76 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
